# Populate the new "Curr-Week" (H) column of the Povit pivot sheet with this
# week's Close/Low/High inputs, fill the same formulas that already exist in
# column G (this week's sibling) across F (last week, previously blank) and
# H (this week, brand new), add the plain Camarilla numbers for H33:H41, and
# move the active-cell selection to H40.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Povit")

# xlPasteFormats - used to copy a cell's formatting (style) without touching
# its contents, mirroring how the "G" column's style is reused for "H" once
# that column starts being populated.
$xlPasteFormats = -4122

function Copy-Format($fromAddr, $toAddr) {
    $ws.Range($fromAddr).Copy() | Out-Null
    $ws.Range($toAddr).PasteSpecial($xlPasteFormats) | Out-Null
}

# ---------------------------------------------------------------------
# Row 2-4: raw inputs (Close / Low / High) for the Pre-Week (F) and
# Curr-Week (H) columns.
# ---------------------------------------------------------------------
$ws.Range("F2").Value = 10540.65
$ws.Range("H2").Value = 10398.35

$ws.Range("F3").Value = 10198.4
$ws.Range("H3").Value = 10198.4

$ws.Range("F4").Value = 10348.049999999999
$ws.Range("H4").Value = 10348.049999999999

# Row 5 separator - F5 no longer carries a (blank) styled cell.
$ws.Range("F5").Clear() | Out-Null

# ---------------------------------------------------------------------
# Pivots block (rows 6-17): F was already styled/blank, just needs the
# formula. H needs both the G-column style and the formula.
# ---------------------------------------------------------------------
$ws.Range("F6").Formula = "=F8+F43"
Copy-Format "G6" "H6"
$ws.Range("H6").Formula = "=H8+H43"

$ws.Range("F7").Formula = "=F11+F43"
Copy-Format "G7" "H7"
$ws.Range("H7").Formula = "=H11+H43"

$ws.Range("F8").Formula = "=(2*F11)-F3"
Copy-Format "G8" "H8"
$ws.Range("H8").Formula = "=(2*H11)-H3"

Copy-Format "G9" "H9"

$ws.Range("F10").Formula = "=F11+F13/2"
Copy-Format "G10" "H10"
$ws.Range("H10").Formula = "=H11+H13/2"

$ws.Range("F11").Formula = "=(F2+F3+F4)/3"
Copy-Format "G11" "H11"
$ws.Range("H11").Formula = "=(H2+H3+H4)/3"

$ws.Range("F12").Formula = "=F11-F13/2"
Copy-Format "G12" "H12"
$ws.Range("H12").Formula = "=H11-H13/2"

$ws.Range("F13").Formula = "=ABS((F11-F46)*2)"
Copy-Format "G13" "H13"
$ws.Range("H13").Formula = "=ABS((H11-H46)*2)"

Copy-Format "G14" "H14"

$ws.Range("F15").Formula = "=2*F11-F2"
Copy-Format "G15" "H15"
$ws.Range("H15").Formula = "=2*H11-H2"

$ws.Range("F16").Formula = "=F11-F43"
Copy-Format "G16" "H16"
$ws.Range("H16").Formula = "=H11-H43"

$ws.Range("F17").Formula = "=F15-F43"
Copy-Format "G17" "H17"
$ws.Range("H17").Formula = "=H15-H43"

Copy-Format "G18" "H18"

# ---------------------------------------------------------------------
# Elliott - Fibonacci block (rows 19-31): F gets the formula; H is a
# brand new cell that needs the G-column style plus the formula.
# ---------------------------------------------------------------------
$ws.Range("F19").Formula = "=(F2/F3)*F4"
Copy-Format "G19" "H19"
$ws.Range("H19").Formula = "=(H2/H3)*H4"

$ws.Range("F20").Formula = "=F21+1.168*(F21-F22)"
Copy-Format "G20" "H20"
$ws.Range("H20").Formula = "=H21+1.168*(H21-H22)"

$ws.Range("F21").Formula = "=F4+F44/2"
Copy-Format "G21" "H21"
$ws.Range("H21").Formula = "=H4+H44/2"

$ws.Range("F22").Formula = "=F4+F44/4"
Copy-Format "G22" "H22"
$ws.Range("H22").Formula = "=H4+H44/4"

$ws.Range("F23").Formula = "=F4+F44/6"
Copy-Format "G23" "H23"
$ws.Range("H23").Formula = "=H4+H44/6"

$ws.Range("F24").Formula = "=F4+F44/12"
Copy-Format "G24" "H24"
$ws.Range("H24").Formula = "=H4+H44/12"

$ws.Range("F25").Formula = "=F4"
Copy-Format "G25" "H25"
$ws.Range("H25").Formula = "=H4"

$ws.Range("F26").Formula = "=F4-F44/12"
Copy-Format "G26" "H26"
$ws.Range("H26").Formula = "=H4-H44/12"

$ws.Range("F27").Formula = "=F4-F44/6"
Copy-Format "G27" "H27"
$ws.Range("H27").Formula = "=H4-H44/6"

$ws.Range("F28").Formula = "=F4-F44/4"
Copy-Format "G28" "H28"
$ws.Range("H28").Formula = "=H4-H44/4"

$ws.Range("F29").Formula = "=F4-F44/2"
Copy-Format "G29" "H29"
$ws.Range("H29").Formula = "=H4-H44/2"

$ws.Range("F30").Formula = "=F29-1.168*(F28-F29)"
Copy-Format "G30" "H30"
$ws.Range("H30").Formula = "=H29-1.168*(H28-H29)"

$ws.Range("F31").Formula = "=F4-(F19-F4)"
Copy-Format "G31" "H31"
$ws.Range("H31").Formula = "=H4-(H19-H4)"

Copy-Format "G32" "H32"

# ---------------------------------------------------------------------
# Camarilla Pivots block (rows 33-41): only H gains cells here (plain
# numbers, mirroring G's style); E/F stay untouched/blank.
# ---------------------------------------------------------------------
Copy-Format "G33" "H33"

Copy-Format "G34" "H34"
$ws.Range("H34").Value = 10559

Copy-Format "G35" "H35"
$ws.Range("H35").Value = 10448

Copy-Format "G36" "H36"
$ws.Range("H36").Value = 10421

# Row 37 ("Close:") keeps its own existing style; only the stray I/J/K
# copies of the formula are dropped now that H37 is meaningful.
$ws.Range("I37:K37").Clear() | Out-Null
$ws.Range("H37").Formula = "=H4"

Copy-Format "G38" "H38"
$ws.Range("H38").Value = 10335

Copy-Format "G39" "H39"
$ws.Range("H39").Value = 10309

Copy-Format "G40" "H40"

Copy-Format "G41" "H41"

Copy-Format "G42" "H42"

# ---------------------------------------------------------------------
# Footer block (rows 43-46): F already styled/blank, needs formula; H
# already styled for rows 43/45 and needs a style refresh (from G) on 46.
# ---------------------------------------------------------------------
$ws.Range("F43").Formula = "=ABS(F2-F3)"
$ws.Range("H43").Formula = "=ABS(H2-H3)"

$ws.Range("F44").Formula = "=F43*1.1"
Copy-Format "G44" "H44"
$ws.Range("H44").Formula = "=H43*1.1"

$ws.Range("F45").Formula = "=(F2+F3)"
$ws.Range("H45").Formula = "=(H2+H3)"

$ws.Range("F46").Formula = "=(F2+F3)/2"
Copy-Format "G46" "H46"
$ws.Range("H46").Formula = "=(H2+H3)/2"

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Move the active selection to H40, as in the final saved state.
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("H40").Select() | Out-Null
